$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '27.131.55'
$ws.Cells.Item(2, 5).Value = '  -0.20%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.626.58'
$ws.Cells.Item(3, 5).Value = '  -0.98%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '214.39'
$ws.Cells.Item(5, 5).Value = '  -1.30%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.523'
$ws.Cells.Item(6, 5).Value = '  +1.85%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.252'
$ws.Cells.Item(8, 5).Value = '  -1.23%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.08%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '20.39'
$ws.Cells.Item(10, 5).Value = '  +1.68%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0847'
$ws.Cells.Item(11, 5).Value = '  -0.03%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.617.08'
$ws.Cells.Item(12, 5).Value = '  -1.64%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.38%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.544'
$ws.Cells.Item(14, 5).Value = '  -0.06%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '27.133.63'
$ws.Cells.Item(15, 5).Value = '  -0.12%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.72'
$ws.Cells.Item(16, 5).Value = '  -3.94%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.62%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '217.02'
$ws.Cells.Item(18, 5).Value = '  -1.05%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.05%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +1.09%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.38'
$ws.Cells.Item(21, 5).Value = '  -0.85%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.43'
$ws.Cells.Item(22, 5).Value = '  -5.59%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.06'
$ws.Cells.Item(23, 5).Value = '  -1.74%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '148.11'
$ws.Cells.Item(24, 5).Value = '  +0.32%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.06%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.31'
$ws.Cells.Item(26, 5).Value = '  -2.99%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.38%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.93%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.42%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.10%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -1.25%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.346.67'
$ws.Cells.Item(33, 5).Value = '  +5.57%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.18%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.40%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0177'
$ws.Cells.Item(36, 5).Value = '  -0.21%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.552'
$ws.Cells.Item(37, 5).Value = '  +1.33%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.28%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.06%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.804'
$ws.Cells.Item(40, 5).Value = '  -0.42%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.23%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +5.84%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(43, 4).Value = '1.765.07'
$ws.Cells.Item(43, 5).Value = '  -1.03%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.23'
$ws.Cells.Item(44, 5).Value = '  -1.49%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '90.72'
$ws.Cells.Item(45, 5).Value = '  -1.24%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.74%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.851'
$ws.Cells.Item(47, 5).Value = '  +28.18%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Cronos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0514'
$ws.Cells.Item(48, 5).Value = '  +0.09%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0998'
$ws.Cells.Item(49, 5).Value = '  +2.24%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.62'
$ws.Cells.Item(50, 5).Value = '  -0.34%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.00'
$ws.Cells.Item(51, 5).Value = '  -0.05%  '
